$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.25'
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '25.55'
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.132'
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05583'
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.483'
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.018'
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8178'
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8411'
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.009713'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9OneONEBestin24h'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1333'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10WazirXWRX'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02864'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11BitrueCoinBTR'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09383'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitMartTokenBMX'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001520'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitForexTokenBF'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006213'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14TigerCashTCH'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.528'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15LEOLEO'

# Row 17
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.022'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16BTSETokenBTSE'

# Row 18
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.3179'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17BitpandaEcosystemTokenBEST'

# Row 19
$ws.Range("B19").Value = 'MandalaExchangeToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06960'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '18MandalaExchangeTokenMDX'

# Row 20
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03210'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '19LiechtensteinCryptoassetsExchangeLCX'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.742'
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04693'
$ws.Range("D23").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004608'
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009699'
$ws.Range("D27").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03654'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006120'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40KickTokenKICK'

# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1052'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41BKEXTokenBKK'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007645'
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005310'
$ws.Range("D45").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002124'
$ws.Range("D48").Style = "Normal"
